$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) Remove the "inventory" sheet - its data is merged into "items"
#    as the new "currentInventory" column.
# -----------------------------------------------------------------
$wb.Worksheets.Item("inventory").Delete() | Out-Null

# -----------------------------------------------------------------
# 2) Update the "items" sheet: rework column D (levelLoad -> strategy)
#    and add currentInventory / FOQ / palletQty columns.
# -----------------------------------------------------------------
$items = $wb.Worksheets.Item("items")

# Header row
$items.Range("D1").Value = "strategy"
$items.Range("E1").Value = "currentInventory"
$items.Range("F1").Value = "FOQ"
$items.Range("G1").Value = "palletQty"

# Column D: boolean levelLoad -> text strategy
$items.Range("D2").Value = "Reactive Level Load - Fast"
$items.Range("D3").Value = "Reactive Level Load - Fast"
$items.Range("D4").Value = "Reactive Level Load - Slow"
$items.Range("D5").Value = "Reactive Level Load - Fast"
$items.Range("D6").Value = "Reactive Level Load - Fast"
$items.Range("D7").Value = "Reactive Level Load - Fast"
$items.Range("D8").Value = "Make-to-Order"

# Column E: currentInventory (values that used to live on the "inventory" sheet)
$items.Range("E2").Value = 33760
$items.Range("E3").Value = 18539
$items.Range("E4").Value = 2929
$items.Range("E5").Value = 12550
$items.Range("E6").Value = 15032
$items.Range("E7").Value = 2986
$items.Range("E8").Value = 1008
$items.Range("E2:E8").NumberFormat = "#,##0"

# Column F: FOQ
$items.Range("F2").Value = 2772
$items.Range("F3").Value = 2772
$items.Range("F4").Value = 2520
$items.Range("F5").Value = 2772
$items.Range("F6").Value = 2772
$items.Range("F7").Value = 2772
$items.Range("F8").Value = 2520

# Column G: palletQty
$items.Range("G2").Value = 504
$items.Range("G3").Value = 504
$items.Range("G4").Value = 504
$items.Range("G5").Value = 504
$items.Range("G6").Value = 504
$items.Range("G7").Value = 504
$items.Range("G8").Value = 504

# Column widths
$items.Columns.Item(4).ColumnWidth = 27.45
$items.Columns.Item(5).ColumnWidth = 8.73

# -----------------------------------------------------------------
# 3) Add the new "calendar" sheet at the end of the workbook.
# -----------------------------------------------------------------
$calendar = $wb.Worksheets.Add([System.Type]::Missing, $wb.Worksheets.Item($wb.Worksheets.Count))
$calendar.Name = "calendar"

$calendar.Range("A1").Value = "period"
$calendar.Range("B1").Value = "startDate"
$calendar.Range("C1").Value = "days"

$calData = @(
    @(202308, 45139, 19),
    @(202309, 45170, 15),
    @(202310, 45200, 18),
    @(202311, 45231, 17),
    @(202312, 45261, 14),
    @(202401, 45292, 18),
    @(202402, 45323, 17),
    @(202403, 45352, 16),
    @(202404, 45383, 17),
    @(202405, 45413, 17),
    @(202406, 45444, 16),
    @(202407, 45474, 18)
)

$r = 2
foreach ($row in $calData) {
    $calendar.Range("A$r").Value = $row[0]
    $calendar.Range("B$r").Value = $row[1]
    $calendar.Range("B$r").NumberFormat = "m/d/yyyy"
    $calendar.Range("C$r").Value = $row[2]
    $r++
}

$calendar.Columns.Item(2).ColumnWidth = 9.45

# -----------------------------------------------------------------
# 4) Add the "days" formula column (=C/20) to the "constraints" sheet.
# -----------------------------------------------------------------
$constraints = $wb.Worksheets.Item("constraints")
$constraints.Range("D1").Value = "days"
$constraints.Range("D2").Formula = "=C2/20"
$constraints.Range("D3:D13").Formula = "=C3/20"

# -----------------------------------------------------------------
# 5) Restore view state: selections on each sheet, and make "items"
#    the active sheet/tab again.
# -----------------------------------------------------------------
$resources = $wb.Worksheets.Item("resources")
$forecast = $wb.Worksheets.Item("forecast")
$requirements = $wb.Worksheets.Item("requirements")

$resources.Activate()
$resources.Range("A13").Select() | Out-Null

$forecast.Activate()
$forecast.Range("B16").Select() | Out-Null

$requirements.Activate()
$requirements.Range("C27").Select() | Out-Null

$constraints.Activate()
$constraints.Range("D2:D13").Select() | Out-Null

$calendar.Activate()
$calendar.Range("C2").Select() | Out-Null

$items.Activate()
$items.Range("D8").Select() | Out-Null
